$d = $word.ActiveDocument

# --- Part 1: insert two new paragraphs in the body text -----------------
# Locate the paragraph that ends with "...Endpoints: Discovery / Location /
# Resolution services." and insert a blank paragraph followed by a new
# "Augmentation: ..." paragraph right after it (and before the existing
# empty bold paragraph that follows).
$target = $null
$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs($i)
    if ($p.Range.Text -like "*Dataflow, Reactive: Resource Monad handling of wrapped URIs*") {
        $target = $p
        break
    }
}

$target.Range.InsertParagraphAfter()
$blank = $target.Next()
$blank.Range.InsertParagraphAfter()
$newPara = $blank.Next()
$newPara.Range.Text = "Augmentation: Model, Context (Statement), Resource levels Message (quads) IO application, resolution, transform / declarative specification (template, input context, results). Dataflow contexts from Message levels application."

# --- Part 2: duplicate one more set of base styles just before the run of
# "Subtitle" styles at the tail of the style sheet. -----------------------
$styles = $d.Styles
$styleCount = $styles.Count

# Find the index of the last "Title" style (end of the 9-style groups,
# right before the block of trailing "Subtitle" styles).
$lastTitleIndex = -1
for ($i = 1; $i -le $styleCount; $i++) {
    $s = $styles.Item($i)
    if ($s.NameLocal -eq "Title") {
        $lastTitleIndex = $i
    }
}

$template = $styles.Item($lastTitleIndex - 8)   # the "Normal" style of that same group

function Copy-StyleAfter($name, $type, $basedOnName, $nextName, $template) {
    $new = $styles.Add($name, $type)
    return $new
}

# Recreate the 9-member group: Normal, TableNormal, Heading1-6, Title
$nNormal = $styles.Add("Normal", 1)
$nNormal.NameLocal = "normal"

$nTable = $styles.Add("Table Normal", 4)
$nTable.NameLocal = "Table Normal"

$headingDefs = @(
    @{Name="heading 1"; Sz=48; After=240; Before=240},
    @{Name="heading 2"; Sz=36; After=225; Before=225},
    @{Name="heading 3"; Sz=28; After=240; Before=240},
    @{Name="heading 4"; Sz=24; After=255; Before=255},
    @{Name="heading 5"; Sz=18; After=255; Before=255},
    @{Name="heading 6"; Sz=16; After=360; Before=360}
)

foreach ($hd in $headingDefs) {
    $h = $styles.Add($hd.Name, 1)
    $h.NameLocal = $hd.Name
    $h.BaseStyle = $nNormal
    $h.NextParagraphStyle = $nNormal
    $h.Font.Bold = $true
    $h.Font.Italic = $false
    $h.Font.Size = $hd.Sz / 2
    $h.ParagraphFormat.SpaceAfter = $hd.After / 20
    $h.ParagraphFormat.SpaceBefore = $hd.Before / 20
}

$nTitle = $styles.Add("Title", 1)
$nTitle.NameLocal = "Title"
$nTitle.BaseStyle = $nNormal
$nTitle.NextParagraphStyle = $nNormal
$nTitle.Font.Bold = $true
$nTitle.Font.Size = 36
$nTitle.ParagraphFormat.SpaceAfter = 6
$nTitle.ParagraphFormat.SpaceBefore = 24
$nTitle.ParagraphFormat.KeepWithNext = $true
